# Auto-generated edit script applying the target diff to before.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: simple single-cell value changes (text swaps / numeric corrections) ---
$cellUpdates = [ordered]@{
  "L7" = 33
  "G8" = "dnasr281@gmail.com, System"
  "L8" = 42
  "G9" = "dnasr281@gmail.com, System"
  "G10" = "dnasr281@gmail.com, System"
  "G12" = "dnasr281@gmail.com, System"
  "G14" = "dnasr281@gmail.com, System"
  "G15" = "dnasr281@gmail.com, System"
  "G17" = "dnasr281@gmail.com, System"
  "P21" = 4
  "Q21" = 4
  "P22" = 4
  "Q22" = 4
  "P23" = 4
  "Q23" = 4
  "P24" = 5
  "Q24" = 4
  "P25" = 4
  "Q25" = 4
  "P26" = 4
  "Q26" = 4
  "G34" = "dnasr281@gmail.com, System"
  "G35" = "dnasr281@gmail.com, System"
  "G36" = "dnasr281@gmail.com, System"
  "G38" = "dnasr281@gmail.com, System"
  "G40" = "dnasr281@gmail.com, System"
  "G41" = "dnasr281@gmail.com, System"
  "G43" = "dnasr281@gmail.com, System"
  "G60" = "dnasr281@gmail.com, System"
  "G61" = "dnasr281@gmail.com, System"
  "G62" = "dnasr281@gmail.com, System"
  "G64" = "dnasr281@gmail.com, System"
  "G66" = "dnasr281@gmail.com, System"
  "G67" = "dnasr281@gmail.com, System"
  "G69" = "dnasr281@gmail.com, System"
  "G86" = "dnasr281@gmail.com, System"
  "G87" = "dnasr281@gmail.com, System"
  "G88" = "dnasr281@gmail.com, System"
  "G90" = "dnasr281@gmail.com, System"
  "G92" = "dnasr281@gmail.com, System"
  "G93" = "dnasr281@gmail.com, System"
  "G95" = "dnasr281@gmail.com, System"
  "G112" = "dnasr281@gmail.com, System"
  "G113" = "dnasr281@gmail.com, System"
  "G114" = "dnasr281@gmail.com, System"
  "G116" = "dnasr281@gmail.com, System"
  "G118" = "dnasr281@gmail.com, System"
  "G119" = "dnasr281@gmail.com, System"
  "G121" = "dnasr281@gmail.com, System"
  "G138" = "dnasr281@gmail.com, System"
  "G139" = "dnasr281@gmail.com, System"
  "G140" = "dnasr281@gmail.com, System"
  "G142" = "dnasr281@gmail.com, System"
  "G144" = "dnasr281@gmail.com, System"
  "G145" = "dnasr281@gmail.com, System"
  "G147" = "dnasr281@gmail.com, System"
}

foreach ($cellRef in $cellUpdates.Keys) {
  $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
}

# --- Part 2: re-style + re-status the "today" session rows (first Pending row per group) ---
# These rows change fill style from the "Pending" (light-yellow, style index 6)
# look to the "Not Recorded" (pink, style index 4) look, and their Status cell text
# flips from "Pending" to "Not Recorded". We copy the *format* from an existing
# style-4 row (A21:I21) onto each target row so the workbook's existing style
# index is reused rather than a new one being minted, then set the Status text.
$formatDonor = $ws.Range("A21:I21")
$formatDonor.Copy()

$notRecordedRows = @(180, 207, 234, 261, 288, 315)
foreach ($r in $notRecordedRows) {
  $target = $ws.Range("A" + $r + ":I" + $r)
  $target.PasteSpecial(-4122)
  $ws.Range("I" + $r).Value = "Not Recorded"
}

$excel.CutCopyMode = 0
